$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G = "Recorded By" holds a comma-separated list of editors, e.g.
# "dnasr281@gmail.com, System". The first two entries in that list need to
# swap places (any trailing entries, if present, stay where they are),
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com" and
# "System, backup@backdoor.com, system" -> "backup@backdoor.com, System, system".

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value()
    if ($text -ne $null -and $text -like "*,*") {
        $parts = $text -split ", "
        if ($parts.Length -ge 2) {
            $first = $parts[0]
            $second = $parts[1]
            $parts[0] = $second
            $parts[1] = $first
            $newText = $parts -join ", "
            $cell.Value = $newText
        }
    }
}
